$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---- Sheet 1: quality_comparison ----
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1: add top+bottom border (style reset to Normal first so no font/alignment carries over)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# D1: add top+bottom+right border (top, right, bottom order keeps every
# intermediate combination matching a pre-existing border definition so no
# throwaway <border> entries get appended to the stylesheet)
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# C2: anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# ---- Sheet 2: computational_comparison ----
$ws2 = $wb.Worksheets.Item("computational_comparison")

# C1 / F1 reuse the same top+bottom style created above (copy formats only,
# avoids the engine allocating throwaway intermediate border/style entries)
$c1.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

# D1 / G1 reuse the same top+bottom+right style created above
$d1.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

# C2, F2: anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5: remove the stray empty inline-string cell entirely
$ws2.Range("G5").ClearContents()
